# Revert "Update mail, registration link and list of committee members of BBS SR"
#
# The roster table (Power-Query table "roster_romandie") lives on the
# worksheet named "Sheet2" (the worksheet physically stored as
# xl/worksheets/sheet1.xml). The previous commit had replaced
# "Marco Eigenmann" and "Dea Putri" with "Laura-Florina Krattinger"; this
# reverts that, restoring Eigenmann/Putri and dropping Krattinger, which
# grows the table from 14 to 15 data rows (A1:E15 -> A1:E16), re-sorted by
# Last name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Grow the query table / range to hold the extra row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E16"))

# Keep the workbook-level external-data named range in sync with the table.
$extData = $wb.Names.Item("ExternalData_1")
$extData.RefersTo = "=Sheet2!`$A`$1:`$E`$16"

# Re-write the roster, sorted alphabetically by Last name (matches the
# table's existing sort state / autofilter).
$ws.Range("A2").Value = "Abdallah"
$ws.Range("B2").Value = "Abouihia"
$ws.Range("C2").Value = "BMS"
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = "Marisa"
$ws.Range("B3").Value = "Bacchi"
$ws.Range("C3").Value = "Statistician/Consultant"
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = "Trista"
$ws.Range("B4").Value = "Bintoro"
$ws.Range("C4").Value = "PMI"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = "Flaminia"
$ws.Range("B5").Value = "Chiesa"
$ws.Range("C5").Value = "Cytel"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = 1

$ws.Range("A6").Value = "Pierre"
$ws.Range("B6").Value = "Colin"
$ws.Range("C6").Value = "BMS"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = 1

$ws.Range("A7").Value = "Liliana"
$ws.Range("B7").Value = "Demenescu"
$ws.Range("C7").Value = "PMI"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 1

$ws.Range("A8").Value = "Elise"
$ws.Range("B8").Value = "Dupuis Lozeron"
$ws.Range("C8").Value = "PMI"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = 1

$ws.Range("A9").Value = "Marco"
$ws.Range("B9").Value = "Eigenmann"
$ws.Range("C9").Value = "PMI"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = 1

$ws.Range("A10").Value = "Alexandra"
$ws.Range("B10").Value = "Korneliou"
$ws.Range("C10").Value = "Biostatistician"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = 1

$ws.Range("A11").Value = "Roland"
$ws.Range("B11").Value = "Marion-Gallois"
$ws.Range("C11").Value = "BMS"
$ws.Range("E11").Value = 1

$ws.Range("A12").Value = "Valeria"
$ws.Range("B12").Value = "Mazzanti"
$ws.Range("C12").Value = "Cytel"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = 1

$ws.Range("A13").Value = "Sandrine"
$ws.Range("B13").Value = "Micallef"
$ws.Range("C13").Value = "Debiopharm"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 1

$ws.Range("A14").Value = "Andrea"
$ws.Range("B14").Value = "Phillips-Beyer"
$ws.Range("C14").Value = "Merck KGaA"
$ws.Range("E14").Value = 1

$ws.Range("A15").Value = "Alessandro"
$ws.Range("B15").Value = "Previtali"
$ws.Range("C15").Value = "BMS"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = 1

$ws.Range("A16").Value = "Dea"
$ws.Range("B16").Value = "Putri"
$ws.Range("C16").Value = "PMI"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = 1

# The old N11:N12 helper cells (beyond the table) are gone in the reverted
# sheet (dimension shrinks from A1:N15 to A1:E16) - clear them out.
$ws.Range("N11:N12").Clear() | Out-Null

# Restore the cursor position recorded in the reverted sheet.
[void]$ws.Range("G7").Select()
